$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# New row 10: Hartley transform sample value + magnitude
$ws.Range("A10").Value = "107839 x10^006 "
$ws.Range("B10").Value = 65536

# Column A width (to fit the new, longer label "107839 x10^006 ")
# (engine quantizes ColumnWidth to the nearest 1/6 of a character; this
# input lands on the stored width closest to the authored 30.7109375)
$ws.Columns.Item(1).ColumnWidth = 29.8333333333333

# Update the selected cell shown in the saved view
$ws.Range("B12").Select()
